# Fix typo in Excel tab name: 'Fundemental data' -> 'Fundamental data'
# plus related selection/data touch-ups captured in the source commit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1) Rename the sheet tab (the actual typo fix).
$ws1.Name = "Fundamental data"

# 2) Populate the missing "Intensity_metric" (Revenue) cells on the
#    "Target data" sheet for the rows using an Intensity target type.
$ws2.Range("D4").Value = "Revenue"
$ws2.Range("D5").Value = "Revenue"
$ws2.Range("D6").Value = "Revenue"
$ws2.Range("D7").Value = "Revenue"
$ws2.Range("D8").Value = "Revenue"
$ws2.Range("D10").Value = "Revenue"
$ws2.Range("D11").Value = "Revenue"
$ws2.Range("D12").Value = "Revenue"
$ws2.Range("D77").Value = "Revenue"
$ws2.Range("D78").Value = "Revenue"

# 3) Restore the selections on each sheet view (and make sure the
#    "Fundamental data" tab stays the active/selected one).
$ws2.Activate() | Out-Null
$ws2.Range("D24").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("I22").Select() | Out-Null
